# Reshape the sheet: shift the existing 2-column (mIoU / Accuracy) table for
# the Open Earth Map dataset down by two rows and right by one column, to make
# room for a new row of merged dataset-title headers plus a "model" column
# header, then add two more (partially filled) dataset blocks to the right
# -- LoveDA dataset Urban and LoveDA dataset Rural -- mirroring the layout of
# the first block but only with the model names copied over so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: 2 new rows at the top, 1 new column on the left. This carries
# the existing data/styles (including the green highlight on the two best
# rows) along with it.
$ws.Rows("1:2").Insert() | Out-Null
$ws.Columns("A:A").Insert() | Out-Null

# Merged, dark-filled, white-centered dataset-title headers on row 2. These
# are entered first so the new shared strings land in the same order as the
# reference workbook (title strings before the "model" header string).
$h1 = $ws.Range("B2:D2")
$h1.Value = "Open Earth Map dataset"
$h1.Font.Color = 16777215
$h1.Interior.Color = 5855577
$h1.HorizontalAlignment = -4108
$h1.Merge() | Out-Null

# Column header for the model-name column (only existing for the first
# block before this edit; the other two blocks are brand new).
$ws.Range("B3").Value = "model"

$h2 = $ws.Range("F2:H2")
$h2.Value = "LoveDA dataset Urban"
$h2.Font.Color = 16777215
$h2.Interior.Color = 5855577
$h2.HorizontalAlignment = -4108
$h2.Merge() | Out-Null

$h3 = $ws.Range("J2:L2")
$h3.Value = "LoveDA dataset Rural"
$h3.Font.Color = 16777215
$h3.Interior.Color = 5855577
$h3.HorizontalAlignment = -4108
$h3.Merge() | Out-Null

$ws.Range("F3").Value = "model"
$ws.Range("G3").Value = "mIoU"
$ws.Range("H3").Value = "Accuracy"

$ws.Range("J3").Value = "model"
$ws.Range("K3").Value = "mIoU"
$ws.Range("L3").Value = "Accuracy"

# The two new blocks only have the model names filled in so far (no
# mIoU/Accuracy numbers yet).
$ws.Range("F4").Value = "U-Net40"
$ws.Range("F5").Value = "SegFormer35"
$ws.Range("F6").Value = "SegFormer30"

$ws.Range("J4").Value = "U-Net40"
$ws.Range("J5").Value = "SegFormer35"
$ws.Range("J6").Value = "SegFormer30"

$ws.Range("J9").Select() | Out-Null
